# Updates cryptos list price/volume figures (and the Maker/VeChain row swap)
# to match the refreshed coinranking.com snapshot described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.155.82"
$ws.Range("E2").Value = "  -6.78%  "
$ws.Range("D3").Value = "1.672.02"
$ws.Range("E3").Value = "  -4.34%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  +0.41%  "
$ws.Range("D5").Value = "'218.61"
$ws.Range("E5").Value = "  -3.29%  "
$ws.Range("E6").Value = "  -12.68%  "
$ws.Range("E7").Value = "  +0.46%  "
$ws.Range("D8").Value = "'0.2635"
$ws.Range("E8").Value = "  -2.78%  "
$ws.Range("D9").Value = "'0.06319"
$ws.Range("E9").Value = "  -4.16%  "
$ws.Range("D10").Value = "'21.60"
$ws.Range("E10").Value = "  -6.85%  "
$ws.Range("D11").Value = "'0.07392"
$ws.Range("E11").Value = "  -1.54%  "
$ws.Range("D12").Value = "1.681.78"
$ws.Range("E12").Value = "  -3.77%  "
$ws.Range("D13").Value = "'4.542"
$ws.Range("E13").Value = "  -3.62%  "
$ws.Range("D14").Value = "'0.5756"
$ws.Range("E14").Value = "  -4.54%  "
$ws.Range("D15").Value = "1.897.89"
$ws.Range("E15").Value = "  -4.38%  "
$ws.Range("D16").Value = "'0.000008515"
$ws.Range("E16").Value = "  -1.27%  "
$ws.Range("D17").Value = "'64.72"
$ws.Range("E17").Value = "  -12.46%  "
$ws.Range("D18").Value = "26.229.36"
$ws.Range("E18").Value = "  -6.49%  "
$ws.Range("E19").Value = "  -6.73%  "
$ws.Range("D20").Value = "'1.006"
$ws.Range("E20").Value = "  +0.21%  "
$ws.Range("D21").Value = "'10.80"
$ws.Range("E21").Value = "  -4.18%  "
$ws.Range("D22").Value = "'186.96"
$ws.Range("E22").Value = "  -8.75%  "
$ws.Range("D23").Value = "'6.178"
$ws.Range("E23").Value = "  -7.15%  "
$ws.Range("E24").Value = "  +0.36%  "
$ws.Range("D25").Value = "'143.09"
$ws.Range("E25").Value = "  -4.89%  "
$ws.Range("D26").Value = "'7.618"
$ws.Range("E26").Value = "  -4.88%  "
$ws.Range("D27").Value = "'0.1164"
$ws.Range("E27").Value = "  -5.78%  "
$ws.Range("D28").Value = "'15.64"
$ws.Range("E28").Value = "  -2.72%  "
$ws.Range("D29").Value = "'1.307"
$ws.Range("E29").Value = "  -5.80%  "
$ws.Range("D30").Value = "'0.05746"
$ws.Range("E30").Value = "  -5.79%  "
$ws.Range("D31").Value = "'1.331"
$ws.Range("E31").Value = "  -3.93%  "
$ws.Range("D32").Value = "'3.502"
$ws.Range("E32").Value = "  -6.09%  "
$ws.Range("D33").Value = "'3.487"
$ws.Range("E33").Value = "  -6.03%  "
$ws.Range("D34").Value = "'1.670"
$ws.Range("E34").Value = "  -0.42%  "
$ws.Range("E35").Value = "  -3.25%  "
$ws.Range("D36").Value = "'0.5973"
$ws.Range("E36").Value = "  -5.93%  "
$ws.Range("E37").Value = "  -2.59%  "
$ws.Range("D38").Value = "'2.631"
$ws.Range("E38").Value = "  -1.38%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.01599"
$ws.Range("E39").Value = "  -4.72%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "1.086.81"
$ws.Range("E40").Value = "  -3.52%  "
$ws.Range("D41").Value = "'5.904"
$ws.Range("E41").Value = "  -5.92%  "
$ws.Range("D42").Value = "'0.8603"
$ws.Range("E42").Value = "  -0.47%  "
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").Value = "'99.89"
$ws.Range("E44").Value = "  +0.29%  "
$ws.Range("D45").Value = "1.819.93"
$ws.Range("E45").Value = "  -4.14%  "
$ws.Range("D46").Value = "'0.00000000111"
$ws.Range("E46").Value = "  +3.10%  "
$ws.Range("D47").Value = "'56.01"
$ws.Range("E47").Value = "  -5.24%  "
$ws.Range("E48").Value = "  -0.67%  "
$ws.Range("D49").Value = "'8.047"
$ws.Range("E49").Value = "  -2.95%  "
$ws.Range("E50").Value = "  -2.68%  "
$ws.Range("D51").Value = "'0.05204"
$ws.Range("E51").Value = "  -3.58%  "
